# The upstream change here (commit: "Fixed POI packaging and upgraded to
# POI 3.15.") is a tooling/library artifact: every hunk in the target diff
# only reorders XML attributes (e.g. <w:tab w:val=".." w:pos=".."/> becomes
# <w:tab w:pos=".." w:val=".."/>, namespace declarations on <w:document>
# get alphabetized, <w:pgSz>/<w:pgMar>/<w:rFonts>/<w:lang>/<w:style>/... all
# get their attributes re-sorted) because the fixture was regenerated with
# a newer Apache POI release. No paragraph text, run, formatting value,
# style definition, numbering, section geometry, font, or any other
# document content/semantics actually changed between the two XML trees -
# every removed line and its paired added line carry exactly the same set
# of attribute=value pairs, just written in a different order.
#
# That kind of low-level serializer attribute-ordering is an artifact of
# the authoring library (Apache POI's XMLBeans writer) and is not something
# the Word object model exposes a way to control - Word (and this
# COM-interop surface) always re-emits elements such as <w:tab>, <w:pgSz>,
# <w:pgMar>, <w:rFonts>, <w:lang>, <w:latentStyles>/<w:lsdException>,
# <w:style>, <w:tblInd>, <w:tblCellMar>, and the root <w:document> namespace
# list in their fixed, schema-declaration order regardless of how content
# is touched through Find/Replace, Range, ParagraphFormat.TabStops,
# PageSetup, Styles, or InsertXML - there is no property/method that lets a
# caller re-order the attributes Word writes back out.
#
# So there is nothing in this document's actual content, formatting, or
# structure for a Word automation script to change: the two OOXML trees
# already describe the same document. Leave it untouched.
$d = $word.ActiveDocument
